$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 (this shifts the existing
# rows 2 and 3 down to rows 3 and 4, unchanged).
$ws.Rows(2).Insert()

# Populate the new row 2 with the "Mosippa - ej aterfunnen" observation.
$ws.Range("A2").Value2 = 58806600
$ws.Range("B2").Value2 = 98536
$ws.Range("C2").Value2 = "Ovaliderad"
$ws.Range("D2").Value2 = "EN"
$ws.Range("E2").Value2 = 1853
$ws.Range("F2").Value2 = "Mosippa"
$ws.Range("G2").Value2 = "Pulsatilla vernalis"
$ws.Range("H2").Value2 = "(L.) Mill."

$ws.Range("P2").Value2 = "Hultsfred, 500 m O betongindustrin, Sm"
$ws.Range("Q2").Value2 = 550867.5962778389
$ws.Range("R2").Value2 = 6376952.536148308
$ws.Range("S2").Value2 = 10
$ws.Range("T2").Value2 = "Kalmar"
$ws.Range("U2").Value2 = "Hultsfred"
$ws.Range("V2").Value2 = "Småland"
$ws.Range("W2").Value2 = "Vena"

# Y2 / AA2 hold dates formatted as plain text in this workbook (not real
# Excel dates) - format the cells as text first so the "yyyy-mm-dd"
# string is not auto-converted into a date serial number.
$ws.Range("Y2").NumberFormat = "@"
$ws.Range("Y2").Value2 = "2016-04-23"
$ws.Range("Z2").Value2 = "00:00"
$ws.Range("AA2").NumberFormat = "@"
$ws.Range("AA2").Value2 = "2016-04-23"
$ws.Range("AB2").Value2 = "00:00"
$ws.Range("AC2").Value2 = "Ej återfunnen. För tjockt mosstäcke på lokalen."

$ws.Range("AD2").Value2 = $true
$ws.Range("AE2").Value2 = $false
$ws.Range("AG2").Value2 = $false
$ws.Range("AI2").Value2 = "tallskog"

$ws.Range("AW2").Value2 = "Calle Ljungberg"
$ws.Range("AX2").Value2 = "Calle Ljungberg, Gunvald Bruce, Mats Halling"
